$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete entire row 6 (its data no longer exists in the target sheet)
$ws.Rows.Item(6).Delete()

$ws.Range("C2").Value = 'No, Statement 1 cannot be inferred from Statement 2. The second statement talks about a game that emulates a board game scenario, but it does not contain any information on whether each board is made up of cells or not.'
$ws.Range("D2").Value = 'No, Statement 1 cannot be implied from Statement 2. Statement 2 does not provide any information about the boards being made up of cells.'
$ws.Range("E2").Value = 'No, Statement 1 cannot be determined from Statement 2. Statement 2 doesn''t provide any information about the boards being made up of cells.'
$ws.Range("F2").Value = 'No, Statement 1 cannot be derived from Statement 2. The second statement does not provide any information about a board being made up of cells.'
$ws.Range("G2").Value = 'No, Statement 1 cannot logically follow from Statement 2. The second statement does not provide information about each board being made up of cells.'
$ws.Range("H2").Value = 'No, Statement 1 cannot be concluded based on Statement 2. Statement 2 doesn''t provide any information about cells and their connection with board games.'
$ws.Range("I2").Value = 'Yes, Statement 2 does support Statement 1 as it implies that the game, being compared to a board game, likely has a segmented structure which could be considered made up of cells, similar to a traditional board game.'
$ws.Range("C3").Value = 'No, Statement 1 cannot be inferred from Statement 2.'
$ws.Range("E3").Value = 'No, Statement 1 cannot be determined from Statement 2.'
$ws.Range("G3").Value = 'No, Statement 1 cannot logically follow from Statement 2 because Statement 2 does not provide information about whether each board is made up of cells.'
$ws.Range("H3").Value = 'No, Statement 1 cannot be concluded based on Statement 2. Statement 2 does not provide information about boards being made up of cells.'
$ws.Range("C4").Value = 'No, Statement 1 cannot be inferred from Statement 2. Statement 2 mentions a board and cells, but it doesn''t directly say that every board is made up of cells.'
$ws.Range("D4").Value = 'No, Statement 1 cannot be implied from Statement 2.'
$ws.Range("G4").Value = 'No, Statement 1 does not logically follow from Statement 2. Statement 2 provides information about a certain game scenario involving a player, a board, and a planet. Statement 1, however, is a general declaration about all boards being made up of cells and does not necessarily follow from the specific game situation laid out in Statement 2.'
$ws.Range("H4").Value = 'No, Statement 1 cannot be concluded based on Statement 2.'
$ws.Range("I4").Value = 'Yes, statement 2 does support statement 1 as it implies that the board is made up of cells that a player can move between.'
$ws.Range("D5").Value = 'No, Statement 1 cannot be implied from Statement 2. The two statements discuss completely different subjects.'
$ws.Range("E5").Value = 'No, Statement 1 cannot be determined from Statement 2. They are about different topics.'
$ws.Range("F5").Value = 'No, Statement 1 cannot be derived from Statement 2. The two statements discuss different topics and do not share any significant connection.'
$ws.Range("G5").Value = 'No, Statement 1 cannot logically follow from Statement 2 as they talk about completely different features of a game – one is about the composition of the board, and the other about the impact of an incorrect hypothesis on a player''s status in the game. The two statements have no logical connection.'
$ws.Range("H5").Value = 'No, Statement 1 cannot be concluded based on Statement 2. The two statements discuss different topics and do not relate to each other.'
